$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 11:52"

# --- Suiza (row 18): updated stats ---
$ws.Range("E18").Value = 6084
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 1593

# --- Bielorrusia moves up above Emiratos Arabes Unidos (rows 36-37) ---
# Row 36 becomes Bielorrusia with fresh stats; row 37 becomes Emiratos Arabes
# Unidos carrying the stats that used to belong to row 36.
$ws.Range("A36").Value = "Bielorrusia"
$ws.Range("B36").Value = 9590
$ws.Range("C36").Value = 817
$ws.Range("D36").Value = 1573
$ws.Range("E36").Value = 7950
$ws.Range("F36").Value = 92
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 67

$ws.Range("A37").Value = "Emiratos Arabes Unidos"
$ws.Range("B37").Value = 9281
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 1760
$ws.Range("E37").Value = 7457
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 64

# --- Kuwait (row 61): updated stats ---
$ws.Range("D61").Value = 656
$ws.Range("E61").Value = 1943

# --- Kazajistan moves up above Barein/Grecia (rows 62-64) ---
# Row 62 becomes Kazajistan with fresh stats; rows 63-64 shift down,
# carrying the stats that used to belong to rows 62-63 (Barein, Grecia).
$ws.Range("A62").Value = "Kazajistan"
$ws.Range("B62").Value = 2525
$ws.Range("C62").Value = 109
$ws.Range("D62").Value = 629
$ws.Range("E62").Value = 1871
$ws.Range("F62").Value = 31
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 25

$ws.Range("A63").Value = "Barein"
$ws.Range("B63").Value = 2518
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 1113
$ws.Range("E63").Value = 1397
$ws.Range("F63").Value = 1
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 8

$ws.Range("A64").Value = "Grecia"
$ws.Range("B64").Value = 2490
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 577
$ws.Range("E64").Value = 1783
$ws.Range("F64").Value = 48
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 130

# --- Uzbekistan (row 68): updated stats ---
$ws.Range("D68").Value = 679
$ws.Range("E68").Value = 1149

# --- Eslovenia (row 79): updated stats ---
$ws.Range("B79").Value = 1388
$ws.Range("C79").Value = 15
$ws.Range("D79").Value = 219
$ws.Range("E79").Value = 1088
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 81

# --- Hong Kong (row 87): updated stats ---
$ws.Range("B87").Value = 1038
$ws.Range("C87").Value = 2
$ws.Range("D87").Value = 753
$ws.Range("E87").Value = 281
$ws.Range("F87").Value = 4

Write-Output "edits applied"
